$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so values like "134.30" or "0.160"
# are not silently coerced into numbers (which would drop trailing zeros
# and change the underlying cell type away from a string).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '58.924.87'
$ws.Cells.Item(2, 5).Value = '  +1.24%  '
$ws.Cells.Item(3, 4).Value = '2.504.38'
$ws.Cells.Item(3, 5).Value = '  -0.55%  '
$ws.Cells.Item(4, 5).Value = '  +0.16%  '
$ws.Cells.Item(5, 4).Value = '535.89'
$ws.Cells.Item(5, 5).Value = '  +3.19%  '
$ws.Cells.Item(6, 4).Value = '134.30'
$ws.Cells.Item(6, 5).Value = '  +1.76%  '
$ws.Cells.Item(7, 5).Value = '  +0.16%  '
$ws.Cells.Item(8, 5).Value = '  +2.47%  '
$ws.Cells.Item(9, 4).Value = '2.508.55'
$ws.Cells.Item(9, 5).Value = '  -0.29%  '
$ws.Cells.Item(10, 4).Value = '0.0996'
$ws.Cells.Item(10, 5).Value = '  +2.32%  '
$ws.Cells.Item(11, 5).Value = '  -2.75%  '
$ws.Cells.Item(12, 5).Value = '  -0.67%  '
$ws.Cells.Item(13, 5).Value = '  -1.32%  '
$ws.Cells.Item(14, 4).Value = '2.951.25'
$ws.Cells.Item(14, 5).Value = '  -0.01%  '
$ws.Cells.Item(15, 4).Value = '58.738.90'
$ws.Cells.Item(15, 5).Value = '  +0.88%  '
$ws.Cells.Item(16, 4).Value = '22.39'
$ws.Cells.Item(16, 5).Value = '  +0.89%  '
$ws.Cells.Item(17, 5).Value = '  +0.53%  '
$ws.Cells.Item(18, 4).Value = '2.505.10'
$ws.Cells.Item(18, 5).Value = '  +0.03%  '
$ws.Cells.Item(19, 4).Value = '10.66'
$ws.Cells.Item(19, 5).Value = '  -0.73%  '
$ws.Cells.Item(20, 4).Value = '4.26'
$ws.Cells.Item(20, 5).Value = '  +1.71%  '
$ws.Cells.Item(21, 4).Value = '321.38'
$ws.Cells.Item(21, 5).Value = '  -0.65%  '
$ws.Cells.Item(22, 4).Value = '6.25'
$ws.Cells.Item(22, 5).Value = '  +3.59%  '
$ws.Cells.Item(23, 5).Value = '  -0.08%  '
$ws.Cells.Item(24, 4).Value = '65.76'
$ws.Cells.Item(24, 5).Value = '  +3.00%  '
$ws.Cells.Item(25, 5).Value = '  +0.63%  '
$ws.Cells.Item(26, 5).Value = '  +1.80%  '
$ws.Cells.Item(27, 4).Value = '0.160'
$ws.Cells.Item(27, 5).Value = '  -1.23%  '
$ws.Cells.Item(28, 4).Value = '7.48'
$ws.Cells.Item(28, 5).Value = '  +1.45%  '
$ws.Cells.Item(29, 4).Value = '0.0₃0757'
$ws.Cells.Item(29, 5).Value = '  +1.07%  '
$ws.Cells.Item(30, 4).Value = '172.20'
$ws.Cells.Item(30, 5).Value = '  +2.83%  '
$ws.Cells.Item(31, 5).Value = '  +1.85%  '
$ws.Cells.Item(32, 4).Value = '6.28'
$ws.Cells.Item(32, 5).Value = '  +0.24%  '
$ws.Cells.Item(33, 5).Value = '  -0.14%  '
$ws.Cells.Item(34, 5).Value = '  +0.15%  '
$ws.Cells.Item(35, 5).Value = '  -0.04%  '
$ws.Cells.Item(36, 4).Value = '18.10'
$ws.Cells.Item(36, 5).Value = '  +0.36%  '
$ws.Cells.Item(37, 5).Value = '  -3.70%  '
$ws.Cells.Item(38, 4).Value = '3.95'
$ws.Cells.Item(38, 5).Value = '  -0.02%  '
$ws.Cells.Item(39, 5).Value = '  +3.89%  '
$ws.Cells.Item(40, 4).Value = '0.831'
$ws.Cells.Item(40, 5).Value = '  +6.75%  '
$ws.Cells.Item(41, 4).Value = '36.57'
$ws.Cells.Item(41, 5).Value = '  -0.58%  '
$ws.Cells.Item(42, 4).Value = '3.49'
$ws.Cells.Item(42, 5).Value = '  +1.61%  '
$ws.Cells.Item(43, 4).Value = '275.09'
$ws.Cells.Item(43, 5).Value = '  -0.82%  '
$ws.Cells.Item(44, 4).Value = '131.85'
$ws.Cells.Item(44, 5).Value = '  +8.06%  '
$ws.Cells.Item(45, 4).Value = '5.03'
$ws.Cells.Item(45, 5).Value = '  -1.29%  '
$ws.Cells.Item(46, 4).Value = '0.590'
$ws.Cells.Item(46, 5).Value = '  -1.34%  '
$ws.Cells.Item(47, 4).Value = '0.0939'
$ws.Cells.Item(47, 5).Value = '  +1.99%  '
$ws.Cells.Item(48, 4).Value = '0.0511'
$ws.Cells.Item(48, 5).Value = '  +2.36%  '
$ws.Cells.Item(49, 5).Value = '  +2.16%  '
$ws.Cells.Item(50, 4).Value = '16.84'
$ws.Cells.Item(50, 5).Value = '  -0.91%  '
$ws.Cells.Item(51, 4).Value = '1.748.51'
$ws.Cells.Item(51, 5).Value = '  +0.27%  '

# Restore the default (Normal) style on column D so no stray number-format
# style is left attached to the cells themselves.
$priceRange.Style = "Normal"

